$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2. Data reporter section: phone number and website were updated
$ws.Range("B9").Value = "0 (312) 32 55 46"
$ws.Range("B10").Value = "www.stat.gov.kg"

# The sheet's active selection moved to B10
$ws.Range("B10").Select()
